$p = $ppt.ActivePresentation

# Slide 11 ("Position Players" content placeholder): remove the
# "Decay = .01" bullet line, leaving an empty, non-bulleted paragraph
# (as if the text had been selected and deleted in the UI).
$s11 = $p.Slides.Item(11)
$shape11 = $s11.Shapes.Item(2)
$para11 = $shape11.TextFrame.TextRange.Paragraphs(3)
$para11.Text = ""
$para11.ParagraphFormat.Bullet.Type = 0

# Slide 12 ("Pitchers" content placeholder): same edit for the
# "Decay = .001" bullet line.
$s12 = $p.Slides.Item(12)
$shape12 = $s12.Shapes.Item(2)
$para12 = $shape12.TextFrame.TextRange.Paragraphs(3)
$para12.Text = ""
$para12.ParagraphFormat.Bullet.Type = 0
